# Update "Pais" worksheet with refreshed COVID-19 country data and timestamp.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 9 de Abril de 2020 a las 08:22"

# --- Row 17: Brasil (stats refreshed, no reordering) ---
$ws.Range("B17").Value = 16195
$ws.Range("C17").Value = 7
$ws.Range("E17").Value = 15246
$ws.Range("G17").Value = 2
$ws.Range("H17").Value = 822

# --- Row 19: Austria (stats refreshed, no reordering) ---
$ws.Range("B19").Value = 12969
$ws.Range("C19").Value = 27
$ws.Range("E19").Value = 8184

# --- Rows 24-25: Australia overtakes Irlanda ---
$ws.Range("A24").Value = "Australia"
$ws.Range("B24").Value = 6104
$ws.Range("C24").Value = 52
$ws.Range("D24").Value = 2813
$ws.Range("E24").Value = 3240
$ws.Range("F24").Value = 84
$ws.Range("G24").Value = 1
$ws.Range("H24").Value = 51

$ws.Range("A25").Value = "Irlanda"
$ws.Range("B25").Value = 6074
$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 25
$ws.Range("E25").Value = 5814
$ws.Range("F25").Value = 165
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 235

# --- Row 47: Tailandia (stats refreshed, no reordering) ---
$ws.Range("B47").Value = 2423
$ws.Range("C47").Value = 54
$ws.Range("D47").Value = 940
$ws.Range("G47").Value = 2
$ws.Range("H47").Value = 32

# --- Rows 75-76: Kazajistan overtakes Camerun ---
$ws.Range("A75").Value = "Kazajistan"
$ws.Range("B75").Value = 759
$ws.Range("C75").Value = 32
$ws.Range("D75").Value = 54
$ws.Range("E75").Value = 698
$ws.Range("F75").Value = 21
$ws.Range("H75").Value = 7

$ws.Range("A76").Value = "Camerun"
$ws.Range("B76").Value = 730
$ws.Range("D76").Value = 60
$ws.Range("E76").Value = 660
$ws.Range("F76").Value = 0
$ws.Range("H76").Value = 10

# --- Row 81: Bulgaria (stats refreshed, no reordering) ---
$ws.Range("B81").Value = 611
$ws.Range("C81").Value = 18
$ws.Range("D81").Value = 48
$ws.Range("E81").Value = 539

# --- Row 95: Taiwan (stats refreshed, no reordering) ---
$ws.Range("B95").Value = 380
$ws.Range("C95").Value = 1
$ws.Range("D95").Value = 80
$ws.Range("E95").Value = 295

# --- Rows 103-105: Kirguistan jumps ahead of Nigeria and Mauricio ---
$ws.Range("A103").Value = "Kirguistan"
$ws.Range("B103").Value = 280
$ws.Range("C103").Value = 10
$ws.Range("D103").Value = 35
$ws.Range("E103").Value = 241
$ws.Range("F103").Value = 5
$ws.Range("H103").Value = 4

$ws.Range("A104").Value = "Nigeria"
$ws.Range("B104").Value = 276
$ws.Range("D104").Value = 44
$ws.Range("E104").Value = 226
$ws.Range("F104").Value = 2
$ws.Range("H104").Value = 6

$ws.Range("A105").Value = "Mauricio"
$ws.Range("B105").Value = 273
$ws.Range("D105").Value = 19
$ws.Range("E105").Value = 247
$ws.Range("F105").Value = 3
$ws.Range("H105").Value = 7

# --- Row 108: Vietnam (stats refreshed, no reordering) ---
$ws.Range("D108").Value = 128
$ws.Range("E108").Value = 123
